$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2-439). The commit bumps that date forward by one day
# (45181 -> 45182, i.e. 2023-09-12 -> 2023-09-13) for all of them.
$lastRow = 439
$range = $ws.Range("C2:C$lastRow")
$range.Value = 45182

$wb.Save()
